$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 23 / Row 24: coins swap position with updated price/volume values
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D23") "97.57"
$ws.Range("E23").Value = "  -7.12%  "

$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D24") "3.59"
$ws.Range("E24").Value = "  -0.21%  "

# Remaining price / volume updates
Set-TextValue $ws.Range("D2") "71.236.31"
$ws.Range("E2").Value = "  -1.76%  "
Set-TextValue $ws.Range("D3") "3.952.36"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue $ws.Range("D5") "537.46"
$ws.Range("E5").Value = "  +3.14%  "
Set-TextValue $ws.Range("D6") "148.44"
$ws.Range("E6").Value = "  +0.57%  "
Set-TextValue $ws.Range("D7") "3.943.01"
$ws.Range("E7").Value = "  -2.86%  "
Set-TextValue $ws.Range("D8") "0.687"
$ws.Range("E8").Value = "  -5.82%  "
Set-TextValue $ws.Range("D9") "1.00"
$ws.Range("E9").Value = "  -0.04%  "
Set-TextValue $ws.Range("D10") "0.742"
$ws.Range("E10").Value = "  -5.22%  "
$ws.Range("E11").Value = "  -6.06%  "
Set-TextValue $ws.Range("D12") "55.37"
$ws.Range("E12").Value = "  +13.39%  "
Set-TextValue $ws.Range("D13") "0.0000318"
$ws.Range("E13").Value = "  -3.91%  "
Set-TextValue $ws.Range("D14") "10.64"
$ws.Range("E14").Value = "  -4.27%  "
Set-TextValue $ws.Range("D15") "4.569.19"
$ws.Range("E15").Value = "  -2.94%  "
Set-TextValue $ws.Range("D16") "3.951.58"
$ws.Range("E16").Value = "  -3.25%  "
Set-TextValue $ws.Range("D17") "13.98"
$ws.Range("E17").Value = "  -3.30%  "
Set-TextValue $ws.Range("D18") "20.55"
$ws.Range("E18").Value = "  -3.89%  "
$ws.Range("E19").Value = "  -1.75%  "
Set-TextValue $ws.Range("D20") "1.17"
$ws.Range("E20").Value = "  -5.80%  "
Set-TextValue $ws.Range("D21") "71.104.50"
$ws.Range("E21").Value = "  -2.05%  "
Set-TextValue $ws.Range("D22") "424.75"
$ws.Range("E22").Value = "  -5.17%  "
$ws.Range("E25").Value = "  +5.50%  "
Set-TextValue $ws.Range("D26") "14.59"
$ws.Range("E26").Value = "  -3.54%  "
Set-TextValue $ws.Range("D27") "11.34"
$ws.Range("E27").Value = "  -0.27%  "
Set-TextValue $ws.Range("D28") "3.86"
Set-TextValue $ws.Range("D29") "10.78"
$ws.Range("E29").Value = "  -2.85%  "
Set-TextValue $ws.Range("D30") "5.89"
$ws.Range("E30").Value = "  +1.01%  "
Set-TextValue $ws.Range("D31") "36.54"
$ws.Range("E31").Value = "  -4.19%  "
Set-TextValue $ws.Range("D32") "7.73"
$ws.Range("E32").Value = "  +16.65%  "
Set-TextValue $ws.Range("D33") "51.17"
$ws.Range("E33").Value = "  +21.25%  "
$ws.Range("E34").Value = "  +0.80%  "
Set-TextValue $ws.Range("D35") "13.35"
$ws.Range("E35").Value = "  -2.66%  "
Set-TextValue $ws.Range("D36") "676.52"
$ws.Range("E36").Value = "  -0.35%  "
Set-TextValue $ws.Range("D37") "65.68"
$ws.Range("E37").Value = "  -3.51%  "
Set-TextValue $ws.Range("D38") "0.443"
$ws.Range("E38").Value = "  +2.87%  "
Set-TextValue $ws.Range("D39") "0.0₃0815"
$ws.Range("E39").Value = "  -5.40%  "
$ws.Range("E40").Value = "  -2.56%  "
Set-TextValue $ws.Range("D41") "3.38"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E43").Value = "  +0.01%  "
Set-TextValue $ws.Range("D44") "0.0483"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("E45").Value = "  -1.19%  "
Set-TextValue $ws.Range("D46") "10.19"
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("E47").Value = "  -6.27%  "
Set-TextValue $ws.Range("D48") "2.66"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  -3.21%  "
Set-TextValue $ws.Range("D50") "3.01"
$ws.Range("E50").Value = "  -2.30%  "
Set-TextValue $ws.Range("D51") "144.35"
$ws.Range("E51").Value = "  +0.30%  "
